# Add a new "Warrant Cancelled Report" sheet, positioned between the
# existing "Warrant Accepted Report" and "Warrant Rejected Report" sheets.
# The new sheet is built from a copy of "Warrant Accepted Report" (same
# layout/styling), then re-pointed at the Warrant-Cancelled NIEM mapping
# strings and trimmed of the two rows that don't apply to this report.

$wb = $excel.ActiveWorkbook

$accepted = $wb.Worksheets.Item("Warrant Accepted Report")
$rejected = $wb.Worksheets.Item("Warrant Rejected Report")

# Duplicate "Warrant Accepted Report" and drop the copy right after it
# (i.e. right before "Warrant Rejected Report").
$accepted.Copy([System.Reflection.Missing]::Value, $accepted)

$newSheet = $wb.Worksheets.Item($accepted.Index + 1)
$newSheet.Name = "Warrant Cancelled Report"

# Two rows present in the "Accepted" mapping do not apply to the
# "Cancelled" mapping - remove them (bottom-most first so the other
# row number stays valid).
$newSheet.Rows.Item(16).Delete()
$newSheet.Rows.Item(13).Delete()

# Title cell.
$newSheet.Range("B1").Value = "Warrant Cancelled Report"

# NIEM mapping column (F) - same shape as "Warrant Accepted Report" but
# rooted at wcr-doc:WarrantAcceptedReport instead of war-doc:WarrantAcceptedReport.
$newSheet.Range("F4").Value = "wcr-doc:WarrantAcceptedReport/nc:Person[@structures:id=/wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderDesignatedSubject/nc:RoleOfPerson/@structures:ref]/nc:PersonName/nc:PersonGivenName"
$newSheet.Range("F5").Value = "wcr-doc:WarrantAcceptedReport/nc:Person[@structures:id=/wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderDesignatedSubject/nc:RoleOfPerson/@structures:ref]/nc:PersonName/nc:PersonMiddleName"
$newSheet.Range("F6").Value = "wcr-doc:WarrantAcceptedReport/nc:Person[@structures:id=/wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderDesignatedSubject/nc:RoleOfPerson/@structures:ref]/nc:PersonName/nc:PersonSurName"
$newSheet.Range("F7").Value = "wcr-doc:WarrantAcceptedReport/nc:Person[@structures:id=/wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderDesignatedSubject/nc:RoleOfPerson/@structures:ref]/nc:PersonName/nc:PersonNameSuffixText"
$newSheet.Range("F8").Value = "wcr-doc:WarrantAcceptedReport/nc:Person[@structures:id=/wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderDesignatedSubject/nc:RoleOfPerson/@structures:ref]/nc:PersonBirthDate/nc:Date"
$newSheet.Range("F9").Value = "wcr-doc:WarrantAcceptedReport/nc:Person[@structures:id=/wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderDesignatedSubject/nc:RoleOfPerson/@structures:ref]/nc:PersonRaceText"
$newSheet.Range("F10").Value = "wcr-doc:WarrantAcceptedReport/nc:Person[@structures:id=/wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderDesignatedSubject/nc:RoleOfPerson/@structures:ref]/nc:PersonSexText"
$newSheet.Range("F13").Value = "wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderEnforcementAgency/wir-ext:AgencyRecordIdentification/nc:IdentificationID"
$newSheet.Range("F14").Value = "wcr-doc:WarrantAcceptedReport/j:Warrant/j:CourtOrderRequestEntity/nc:EntityPerson/wir-ext:PersonEmployeeIdentification/nc:IdentificationID"
$newSheet.Range("F15").Value = "wcr-doc:WarrantAcceptedReport/j:Warrant/wir-ext:WarrantAugmentation/wir-ext:StateWarrantRepositoryIdentification/nc:IdentificationID/#text"

# Make the new sheet the active tab/selection, matching the saved view
# state of the authored workbook.
$newSheet.Activate()
$newSheet.Range("A2").Select()

$wb.Save()
